# This script updates market-data columns (H-N) on several sheets of the
# Gungnir_Profits workbook to reflect a refreshed data pull from the scheduled
# runner. Each block targets specific leve rows identified by cell address.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31: Hush Little Wailer
$ws.Range("H31").Value = 711.8333
$ws.Range("I31").Value = 711.8333
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2135.4999
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1905.4999
# Row 88: The Grave of Hemlock Groves
$ws.Range("H88").Value = 9128318
$ws.Range("I88").Value = 1165.375
$ws.Range("J88").Value = 15213086
$ws.Range("K88").Value = 1165.375
$ws.Range("L88").Value = 15213086
$ws.Range("M88").Value = -759.375
$ws.Range("N88").Value = -15213898
# Row 91: Dappling the Highlands (L)
$ws.Range("H91").Value = 9128318
$ws.Range("I91").Value = 1165.375
$ws.Range("J91").Value = 15213086
$ws.Range("K91").Value = 1165.375
$ws.Range("L91").Value = 15213086
$ws.Range("M91").Value = 238.625
$ws.Range("N91").Value = -15215894
# Row 106: Making Your Mark
$ws.Range("H106").Value = 90911790
$ws.Range("I106").Value = 200001860
$ws.Range("J106").Value = 3400
$ws.Range("K106").Value = 200001860
$ws.Range("L106").Value = 3400
$ws.Range("M106").Value = -200001229
$ws.Range("N106").Value = -4662
# Row 112: Making Ends Meet
$ws.Range("H112").Value = 12285.028
$ws.Range("I112").Value = 466.66666
$ws.Range("J112").Value = 12806.427
$ws.Range("K112").Value = 1399.99998
$ws.Range("L112").Value = 38419.281
$ws.Range("M112").Value = -291.9999800000001
$ws.Range("N112").Value = -40635.281
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 60615828
$ws.Range("I132").Value = 51731804
$ws.Range("J132").Value = 125025000
$ws.Range("K132").Value = 155195412
$ws.Range("L132").Value = 375075000
$ws.Range("M132").Value = -155192882
$ws.Range("N132").Value = -375080060
# Row 135: For Tired Minds
$ws.Range("H135").Value = 886.283
$ws.Range("I135").Value = 752.5625
$ws.Range("J135").Value = 2170
$ws.Range("K135").Value = 6773.0625
$ws.Range("L135").Value = 19530
$ws.Range("M135").Value = -4238.0625
$ws.Range("N135").Value = -24600

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 1980.28
$ws.Range("I32").Value = 1908.7386
$ws.Range("J32").Value = 2504.9167
$ws.Range("K32").Value = 1908.7386
$ws.Range("L32").Value = 2504.9167
$ws.Range("M32").Value = -1621.7386
$ws.Range("N32").Value = -3078.9167

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 22139.6
$ws.Range("I22").Value = 232.66667
$ws.Range("J22").Value = 55000
$ws.Range("K22").Value = 232.66667
$ws.Range("L22").Value = 55000
$ws.Range("M22").Value = 117.33333
$ws.Range("N22").Value = -55700
# Row 31: Wall Not Found
$ws.Range("H31").Value = 1375.4
$ws.Range("I31").Value = 983.3200000000001
$ws.Range("J31").Value = 1702.1333
$ws.Range("K31").Value = 983.3200000000001
$ws.Range("L31").Value = 1702.1333
$ws.Range("M31").Value = -688.3200000000001
$ws.Range("N31").Value = -2292.1333
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 1375.4
$ws.Range("I34").Value = 983.3200000000001
$ws.Range("J34").Value = 1702.1333
$ws.Range("K34").Value = 983.3200000000001
$ws.Range("L34").Value = 1702.1333
$ws.Range("M34").Value = -781.3200000000001
$ws.Range("N34").Value = -2106.1333
# Row 50: The Arsenal of Theocracy
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").ClearContents()
# Row 51: Greenstone for Greenhorns
$ws.Range("H51").Value = 17599.6
$ws.Range("I51").Value = 9800
$ws.Range("J51").Value = 19549.5
$ws.Range("K51").Value = 9800
$ws.Range("L51").Value = 19549.5
$ws.Range("M51").Value = -9064
$ws.Range("N51").Value = -21021.5
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 11765379
$ws.Range("I58").Value = 19231300
$ws.Range("J58").Value = 898
$ws.Range("K58").Value = 19231300
$ws.Range("L58").Value = 898
$ws.Range("M58").Value = -19231097
$ws.Range("N58").Value = -1304
# Row 61: Incant Now, Think Later
$ws.Range("H61").Value = 17599.6
$ws.Range("I61").Value = 9800
$ws.Range("J61").Value = 19549.5
$ws.Range("K61").Value = 9800
$ws.Range("L61").Value = 19549.5
$ws.Range("M61").Value = -9452
$ws.Range("N61").Value = -20245.5
# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 5051429.5
$ws.Range("I132").Value = 691.3214
$ws.Range("J132").Value = 33335564
$ws.Range("K132").Value = 2073.9642
$ws.Range("L132").Value = 100006692
$ws.Range("M132").Value = 456.0357999999997
$ws.Range("N132").Value = -100011752
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 822.1096
$ws.Range("I134").Value = 685.76666
$ws.Range("J134").Value = 1451.3846
$ws.Range("K134").Value = 2057.29998
$ws.Range("L134").Value = 4354.1538
$ws.Range("M134").Value = 477.7000200000002
$ws.Range("N134").Value = -9424.1538
# Row 136: Turali Quality
$ws.Range("H136").Value = 11765379
$ws.Range("I136").Value = 19231300
$ws.Range("J136").Value = 898
$ws.Range("K136").Value = 57693900
$ws.Range("L136").Value = 2694
$ws.Range("M136").Value = -57691350
$ws.Range("N136").Value = -7794

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 31534792
$ws.Range("I5").Value = 27778138
$ws.Range("J5").Value = 38470150
$ws.Range("K5").Value = 83334414
$ws.Range("L5").Value = 115410450
$ws.Range("M5").Value = -83334302
$ws.Range("N5").Value = -115410674
# Row 56: Culture Club
$ws.Range("H56").Value = 9168.333000000001
$ws.Range("I56").Value = 9168.333000000001
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 9168.333000000001
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -8638.333000000001
# Row 105: Fish Box
$ws.Range("H105").Value = 5176.3335
$ws.Range("I105").Value = 3000
$ws.Range("J105").Value = 5611.6
$ws.Range("K105").Value = 9000
$ws.Range("L105").Value = 16834.8
$ws.Range("M105").Value = -6379
$ws.Range("N105").Value = -22076.8
# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 31534792
$ws.Range("I135").Value = 27778138
$ws.Range("J135").Value = 38470150
$ws.Range("K135").Value = 250003242
$ws.Range("L135").Value = 346231350
$ws.Range("M135").Value = -250000707
$ws.Range("N135").Value = -346236420

$ws = $wb.Worksheets.Item("GSM")
# Row 126: Gold Rush Order
$ws.Range("H126").Value = 2236.1428
$ws.Range("I126").Value = 2006
$ws.Range("J126").Value = 2274.5
$ws.Range("K126").Value = 6018
$ws.Range("L126").Value = 6823.5
$ws.Range("M126").Value = -3548
$ws.Range("N126").Value = -11763.5

$ws = $wb.Worksheets.Item("WVR")
# Row 3: Trew Enough
$ws.Range("H3").Value = 702
$ws.Range("I3").Value = 702
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 702
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -588
$ws.Range("N3").ClearContents()
# Row 7: Long Hair, Long Life
$ws.Range("H7").Value = 931.8
$ws.Range("I7").Value = 614.75
$ws.Range("J7").Value = 2200
$ws.Range("K7").Value = 614.75
$ws.Range("L7").Value = 2200
$ws.Range("M7").Value = -501.75
$ws.Range("N7").Value = -2426
# Row 9: A Taste for Dalmaticae
$ws.Range("H9").Value = 27877.25
$ws.Range("I9").Value = 2599.4
$ws.Range("J9").Value = 70007
$ws.Range("K9").Value = 2599.4
$ws.Range("L9").Value = 70007
$ws.Range("M9").Value = -2459.4
$ws.Range("N9").Value = -70287
# Row 12: This Is Why You Can't Have Nice Things
$ws.Range("H12").Value = 80007
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 80007
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 80007
$ws.Range("N12").Value = -80291
# Row 14: Hat in Hand
$ws.Range("H14").Value = 89900
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 89900
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 89900
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -90236
# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1579.5
$ws.Range("I126").Value = 1297.5
$ws.Range("J126").Value = 1767.5
$ws.Range("K126").Value = 3892.5
$ws.Range("L126").Value = 5302.5
$ws.Range("M126").Value = -1422.5
$ws.Range("N126").Value = -10242.5
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 6344904.5
$ws.Range("I132").Value = 18475.982
$ws.Range("J132").Value = 25007868
$ws.Range("K132").Value = 55427.946
$ws.Range("L132").Value = 75023604
$ws.Range("M132").Value = -52897.946
$ws.Range("N132").Value = -75028664
